$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean slate for the populated area so styles/strings are rebuilt
# in the exact order required by the target layout.
$ws.Range("A1:K5").Clear()

# --- Header row (row 1) ---
# Plain (unstyled) text headers for the new identifying columns.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"

# Unit headers (F1:K1) keep the font-9 styling used throughout the sheet,
# but without a number-format override (General). We build this style via a
# temporary named style so we get a distinct cell format entry, then drop the
# named style itself so only the plain format remains.
$hdrStyle = $wb.Styles.Add("UnitHeaderStyle")
$hdrStyle.Font.Name = "Arial"
$hdrStyle.Font.Size = 9

$ws.Range("F1").Value = "(m3/s)"
$ws.Range("F1").Style = "UnitHeaderStyle"

$ws.Range("G1").Value = "(MW1)"
$ws.Range("G1").Style = "UnitHeaderStyle"

$ws.Range("H1").Value = "(MW2)"
$ws.Range("H1").Style = "UnitHeaderStyle"

$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("I1").Style = "UnitHeaderStyle"

$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("J1").Style = "UnitHeaderStyle"

$ws.Range("K1").Value = "(GWh) Year"
$ws.Range("K1").Style = "UnitHeaderStyle"

$wb.Styles.Item("UnitHeaderStyle").Delete()

# --- Row 2: Kembs ---
$ws.Range("A2").Value = 1
$ws.Range("A2").Font.Name = "Arial"
$ws.Range("A2").Font.Size = 9
$ws.Range("A2").NumberFormat = "0"

$ws.Range("B2").Value = 110500
$ws.Range("B2").Font.Name = "Arial"
$ws.Range("B2").Font.Size = 9
$ws.Range("B2").NumberFormat = "0"

$ws.Range("C2").Value = "Kembs"
$ws.Range("C2").Font.Name = "Arial"
$ws.Range("C2").Font.Size = 9

$ws.Range("D2").Value = 1932
$ws.Range("D2").Font.Name = "Arial"
$ws.Range("D2").Font.Size = 9
$ws.Range("D2").NumberFormat = "0"

$ws.Range("E2").Value = 1983
$ws.Range("E2").Font.Name = "Arial"
$ws.Range("E2").Font.Size = 9
$ws.Range("E2").NumberFormat = "0"

$ws.Range("F2").Value = 1400
$ws.Range("F2").Font.Name = "Arial"
$ws.Range("F2").Font.Size = 9
$ws.Range("F2").NumberFormat = "#,##0.00"

$ws.Range("G2").Value = 31
$ws.Range("G2").Font.Name = "Arial"
$ws.Range("G2").Font.Size = 9
$ws.Range("G2").NumberFormat = "0.00"

$ws.Range("H2").Value = 31.5
$ws.Range("H2").Font.Name = "Arial"
$ws.Range("H2").Font.Size = 9
$ws.Range("H2").NumberFormat = "0.00"

$ws.Range("I2").Value = 75.400000000000006
$ws.Range("I2").Font.Name = "Arial"
$ws.Range("I2").Font.Size = 9
$ws.Range("I2").NumberFormat = "0.00"

$ws.Range("J2").Value = 99.4
$ws.Range("J2").Font.Name = "Arial"
$ws.Range("J2").Font.Size = 9
$ws.Range("J2").NumberFormat = "0.00"

$ws.Range("K2").Value = 174.8
$ws.Range("K2").Font.Name = "Arial"
$ws.Range("K2").Font.Size = 9
$ws.Range("K2").NumberFormat = "0.00"

# --- Row 3: Birsfelden ---
$ws.Range("A3").Value = 2
$ws.Range("A3").Font.Name = "Arial"
$ws.Range("A3").Font.Size = 9
$ws.Range("A3").NumberFormat = "0"

$ws.Range("B3").Value = 109700
$ws.Range("B3").Font.Name = "Arial"
$ws.Range("B3").Font.Size = 9
$ws.Range("B3").NumberFormat = "0"

$ws.Range("C3").Value = "Birsfelden"
$ws.Range("C3").Font.Name = "Arial"
$ws.Range("C3").Font.Size = 9

$ws.Range("D3").Value = 1955
$ws.Range("D3").Font.Name = "Arial"
$ws.Range("D3").Font.Size = 9
$ws.Range("D3").NumberFormat = "0"

$ws.Range("E3").Value = 1999
$ws.Range("E3").Font.Name = "Arial"
$ws.Range("E3").Font.Size = 9
$ws.Range("E3").NumberFormat = "0"

$ws.Range("F3").Value = 1500
$ws.Range("F3").Font.Name = "Arial"
$ws.Range("F3").Font.Size = 9
$ws.Range("F3").NumberFormat = "#,##0.00"

$ws.Range("G3").Value = 15.55
$ws.Range("G3").Font.Name = "Arial"
$ws.Range("G3").Font.Size = 9
$ws.Range("G3").NumberFormat = "0.00"

$ws.Range("H3").Value = 15.16
$ws.Range("H3").Font.Name = "Arial"
$ws.Range("H3").Font.Size = 9
$ws.Range("H3").NumberFormat = "0.00"

$ws.Range("I3").Value = 38.1
$ws.Range("I3").Font.Name = "Arial"
$ws.Range("I3").Font.Size = 9
$ws.Range("I3").NumberFormat = "0.00"

$ws.Range("J3").Value = 48.52
$ws.Range("J3").Font.Name = "Arial"
$ws.Range("J3").Font.Size = 9
$ws.Range("J3").NumberFormat = "0.00"

$ws.Range("K3").Value = 86.61
$ws.Range("K3").Font.Name = "Arial"
$ws.Range("K3").Font.Size = 9
$ws.Range("K3").NumberFormat = "0.00"

# --- Row 4: Kembs-Centrale de dotation ---
$ws.Range("A4").Value = 3
$ws.Range("A4").Font.Name = "Arial"
$ws.Range("A4").Font.Size = 9
$ws.Range("A4").NumberFormat = "0"

$ws.Range("B4").Value = 110490
$ws.Range("B4").Font.Name = "Arial"
$ws.Range("B4").Font.Size = 9
$ws.Range("B4").NumberFormat = "0"

$ws.Range("C4").Value = "Kembs-Centrale de dotation"
$ws.Range("C4").Font.Name = "Arial"
$ws.Range("C4").Font.Size = 9

$ws.Range("D4").Value = 1966
$ws.Range("D4").Font.Name = "Arial"
$ws.Range("D4").Font.Size = 9
$ws.Range("D4").NumberFormat = "0"

# Note: no Date End (E4) for this plant.

$ws.Range("F4").Value = 27
$ws.Range("F4").Font.Name = "Arial"
$ws.Range("F4").Font.Size = 9
$ws.Range("F4").NumberFormat = "0.00"

$ws.Range("G4").Value = 0.56999999999999995
$ws.Range("G4").Font.Name = "Arial"
$ws.Range("G4").Font.Size = 9
$ws.Range("G4").NumberFormat = "0.00"

$ws.Range("H4").Value = 0.63
$ws.Range("H4").Font.Name = "Arial"
$ws.Range("H4").Font.Size = 9
$ws.Range("H4").NumberFormat = "0.00"

$ws.Range("I4").Value = 1.8
$ws.Range("I4").Font.Name = "Arial"
$ws.Range("I4").Font.Size = 9
$ws.Range("I4").NumberFormat = "0.00"

$ws.Range("J4").Value = 2.2000000000000002
$ws.Range("J4").Font.Name = "Arial"
$ws.Range("J4").Font.Size = 9
$ws.Range("J4").NumberFormat = "0.00"

$ws.Range("K4").Value = 4
$ws.Range("K4").Font.Name = "Arial"
$ws.Range("K4").Font.Size = 9
$ws.Range("K4").NumberFormat = "0.00"

# Match the saved selection state from the target workbook.
$ws.Range("A2:K2").Select()
